$d = $word.ActiveDocument

# Locate the paragraph that holds the `print(...)` line (it currently ends
# with a run of otherwise-empty runs). We split it into two paragraphs:
#   1) the print(...) text plus the two runs that immediately follow it
#   2) a brand-new paragraph (same sz/highlight rPr) whose first run now
#      carries 7 literal spaces, followed by one more empty run; the final
#      trailing empty run is dropped entirely.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "print(*") {
        $target = $cand
        break
    }
}

$replacementXml = @"
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:highlight w:val="none"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:highlight w:val="none"/>
              </w:rPr>
              <w:t xml:space="preserve">print(«Привет Мир!»)</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:highlight w:val="none"/>
              </w:rPr>
            </w:r>
            <w:r></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:highlight w:val="none"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:highlight w:val="none"/>
              </w:rPr>
              <w:t xml:space="preserve">       </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:sz w:val="22"/>
                <w:highlight w:val="none"/>
              </w:rPr>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$target.Range.InsertXML($replacementXml)
